$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5849.75
$ws.Range("I43").Value = 4440.4
$ws.Range("K43").Value = 4440.4
$ws.Range("M43").Value = -4371.4

$ws.Range("H86").Value = 2546.6667
$ws.Range("I86").Value = 2399
$ws.Range("K86").Value = 2399
$ws.Range("M86").Value = -1276

$ws.Range("H89").Value = 2546.6667
$ws.Range("I89").Value = 2399
$ws.Range("K89").Value = 11995
$ws.Range("M89").Value = -6379

$ws.Range("H116").Value = 10562.4375
$ws.Range("I116").Value = 10880.637
$ws.Range("J116").Value = 9862.4
$ws.Range("K116").Value = 10880.637
$ws.Range("L116").Value = 9862.4
$ws.Range("M116").Value = -7438.637000000001
$ws.Range("N116").Value = -16746.4

$ws.Range("H132").Value = 3988981.2
$ws.Range("I132").Value = 5743603.5
$ws.Range("K132").Value = 17230810.5
$ws.Range("M132").Value = -17228280.5

$ws.Range("H137").Value = 9876.450000000001
$ws.Range("I137").Value = 13214.777
$ws.Range("K137").Value = 39644.331
$ws.Range("M137").Value = -37094.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18440.047
$ws.Range("I32").Value = 20248.715
$ws.Range("K32").Value = 20248.715
$ws.Range("M32").Value = -19961.715

$ws.Range("H45").Value = 3145.4707
$ws.Range("I45").Value = 2048.0715
$ws.Range("J45").Value = 8266.666999999999
$ws.Range("K45").Value = 2048.0715
$ws.Range("L45").Value = 8266.666999999999
$ws.Range("M45").Value = -1671.0715
$ws.Range("N45").Value = -9020.666999999999

$ws.Range("H61").Value = 10535.272
$ws.Range("I61").Value = 1209.7778
$ws.Range("J61").Value = 52500
$ws.Range("K61").Value = 1209.7778
$ws.Range("L61").Value = 52500
$ws.Range("M61").Value = -997.7778000000001
$ws.Range("N61").Value = -52924

$ws.Range("H74").Value = 1221801.2
$ws.Range("I74").Value = 1502251.5
$ws.Range("K74").Value = 1502251.5
$ws.Range("M74").Value = -1501377.5

$ws.Range("H77").Value = 1221801.2
$ws.Range("I77").Value = 1502251.5
$ws.Range("K77").Value = 7511257.5
$ws.Range("M77").Value = -7506889.5

$ws.Range("H97").Value = 1153.2433
$ws.Range("I97").Value = 916.1852
$ws.Range("J97").Value = 1793.3
$ws.Range("K97").Value = 916.1852
$ws.Range("L97").Value = 1793.3
$ws.Range("M97").Value = -420.1852
$ws.Range("N97").Value = -2785.3

$ws.Range("H102").Value = 3643.1333
$ws.Range("I102").Value = 3760.5
$ws.Range("K102").Value = 3760.5
$ws.Range("M102").Value = -2138.5

$ws.Range("H135").Value = 88995
$ws.Range("J135").Value = 88995
$ws.Range("L135").Value = 88995
$ws.Range("N135").Value = -99135

$ws.Range("H136").Value = 10535.272
$ws.Range("I136").Value = 1209.7778
$ws.Range("J136").Value = 52500
$ws.Range("K136").Value = 3629.3334
$ws.Range("L136").Value = 157500
$ws.Range("M136").Value = -1079.3334
$ws.Range("N136").Value = -162600

$ws.Range("H137").Value = 111719.5
$ws.Range("J137").Value = 119132.78
$ws.Range("L137").Value = 119132.78
$ws.Range("N137").Value = -129332.78

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2042.1305
$ws.Range("I134").Value = 1558.3823
$ws.Range("J134").Value = 3412.75
$ws.Range("K134").Value = 4675.1469
$ws.Range("L134").Value = 10238.25
$ws.Range("M134").Value = -2140.1469
$ws.Range("N134").Value = -15308.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22490.312
$ws.Range("I58").Value = 2391.6365
$ws.Range("K58").Value = 2391.6365
$ws.Range("M58").Value = -2188.6365

$ws.Range("H69").Value = 57500
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61498

$ws.Range("H72").Value = 57500
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -187488

$ws.Range("H99").Value = 12238.9
$ws.Range("I99").Value = 9246.25
$ws.Range("J99").Value = 14234
$ws.Range("K99").Value = 9246.25
$ws.Range("L99").Value = 14234
$ws.Range("M99").Value = -7748.25
$ws.Range("N99").Value = -17230

$ws.Range("H126").Value = 12238.9
$ws.Range("I126").Value = 9246.25
$ws.Range("J126").Value = 14234
$ws.Range("K126").Value = 27738.75
$ws.Range("L126").Value = 42702
$ws.Range("M126").Value = -25268.75
$ws.Range("N126").Value = -47642

$ws.Range("H132").Value = 35495.344
$ws.Range("I132").Value = 36584.465
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 109753.395
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -107223.395
$ws.Range("N132").Value = -20060

$ws.Range("H134").Value = 1075.2142
$ws.Range("I134").Value = 1003.25
$ws.Range("J134").Value = 1507
$ws.Range("K134").Value = 3009.75
$ws.Range("L134").Value = 4521
$ws.Range("M134").Value = -474.75
$ws.Range("N134").Value = -9591

$ws.Range("H135").Value = 119997.336
$ws.Range("J135").Value = 119997.336
$ws.Range("L135").Value = 119997.336
$ws.Range("N135").Value = -130137.336

$ws.Range("H136").Value = 22490.312
$ws.Range("I136").Value = 2391.6365
$ws.Range("K136").Value = 7174.9095
$ws.Range("M136").Value = -4624.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 7142.857
$ws.Range("I63").Value = 50000
$ws.Range("K63").Value = 150000
$ws.Range("M63").Value = -149251

$ws.Range("H66").Value = 7142.857
$ws.Range("I66").Value = 50000
$ws.Range("K66").Value = 450000
$ws.Range("M66").Value = -446256

$ws.Range("H70").Value = 4925
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -10185
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 4925
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -9408
$ws.Range("N73").Value = -17184

$ws.Range("H93").Value = 5500
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744

$ws.Range("H132").Value = 1819.9166
$ws.Range("I132").Value = 2139.2856
$ws.Range("J132").Value = 1372.8
$ws.Range("K132").Value = 19253.5704
$ws.Range("L132").Value = 12355.2
$ws.Range("M132").Value = -16723.5704
$ws.Range("N132").Value = -17415.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 622.55554
$ws.Range("I2").Value = 972.6667
$ws.Range("J2").Value = 272.44446
$ws.Range("K2").Value = 972.6667
$ws.Range("L2").Value = 272.44446
$ws.Range("M2").Value = -859.6667
$ws.Range("N2").Value = -498.44446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1655.5555
$ws.Range("J22").Value = 1757.1428
$ws.Range("L22").Value = 1757.1428
$ws.Range("N22").Value = -2347.1428

$ws.Range("H27").Value = 1655.5555
$ws.Range("J27").Value = 1757.1428
$ws.Range("L27").Value = 1757.1428
$ws.Range("N27").Value = -1971.1428

$ws.Range("H63").Value = 71995
$ws.Range("J63").Value = 71995
$ws.Range("L63").Value = 71995
$ws.Range("N63").Value = -73493

$ws.Range("H66").Value = 71995
$ws.Range("J66").Value = 71995
$ws.Range("L66").Value = 215985
$ws.Range("N66").Value = -223473

$ws.Range("H93").Value = 2337.2727
$ws.Range("I93").Value = 2655.1667
$ws.Range("J93").Value = 1955.8
$ws.Range("K93").Value = 2655.1667
$ws.Range("L93").Value = 1955.8
$ws.Range("M93").Value = -1407.1667
$ws.Range("N93").Value = -4451.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17036.5
$ws.Range("I136").Value = 19120.969
$ws.Range("J136").Value = 3279
$ws.Range("K136").Value = 57362.90700000001
$ws.Range("L136").Value = 9837
$ws.Range("M136").Value = -54812.90700000001
$ws.Range("N136").Value = -14937
